# Rework of section "Figures in Chapter 3" (3.2 - 3.14): each caption shifts
# up to take the place of the one before it (the old "Research Activities and
# Contexts" caption for Figure 3.2 is dropped), the "Ideation Grids" caption is
# renamed to "Ideation Decks" as it moves, and a brand-new final caption is
# appended for what is now the last figure in the list (3.14).
#
# Applied back-to-front so every Find.Execute only ever matches the single,
# still-untouched occurrence of its source text.

$d = $word.ActiveDocument

$ldq = [char]0x201C   # “
$rdq = [char]0x201D   # ”
$en  = [char]0x2013   # –

function Replace-Caption([string]$oldText, [string]$newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
    if (-not $found) {
        throw "Caption not found: $oldText"
    }
}

# 3.14 (was 3.13's text): "- Pilot Study Recruitment Poster" -> new final caption
Replace-Caption "- Pilot Study Recruitment Poster" "- How the Case Studies and Peripheral Activities Contribute to This Thesis"

# 3.13 (was 3.12's text)
Replace-Caption "- Spreadsheet-based Quantitative Analysis of Interview Data for Case Study Two" "- Pilot Study Recruitment Poster"

# 3.12 (was 3.11's text)
Replace-Caption "- Thematic Analysis of Qualitative Data using Quirkos for Case Study One" "- Spreadsheet-based Quantitative Analysis of Interview Data for Case Study Two"

# 3.11 (was 3.10's text)
Replace-Caption "- Storyboarding Cards $en A Collaboratively-constructed Narrative Created through Discussion From a Palette of Possible Parent and Staff Actions" "- Thematic Analysis of Qualitative Data using Quirkos for Case Study One"

# 3.10 (was 3.9's text)
Replace-Caption "- Group Poster Design $en A Participant-designed Poster to Advertise Features of Imagined Data Interface Products" "- Storyboarding Cards $en A Collaboratively-constructed Narrative Created through Discussion From a Palette of Possible Parent and Staff Actions"

# 3.9 (was 3.8's text, renamed "Ideation Grids" -> "Ideation Decks")
Replace-Caption "- Ideation Grids $en Combining Random Design Ingredients to Generate New Ideas" "- Group Poster Design $en A Participant-designed Poster to Advertise Features of Imagined Data Interface Products"

# 3.8 (was 3.7's text)
Replace-Caption "- Home Interviewing: Card Sorting With a Family in Their Living Room" "- Ideation Decks $en Combining Random Design Ingredients to Generate New Ideas"

# 3.7 (was 3.6's text)
Replace-Caption "- Personal Data Examples $en Making Data Relatable" "- Home Interviewing: Card Sorting With a Family in Their Living Room"

# 3.6 (was 3.5's text)
Replace-Caption "- Family Civic Data Cards $en Things to Think With" "- Personal Data Examples $en Making Data Relatable"

# 3.5 (was 3.4's text)
Replace-Caption "- Sentence Ranking $en Bringing Support Workers and Families to a Shared Problem Space" "- Family Civic Data Cards $en Things to Think With"

# 3.4 (was 3.3's text)
Replace-Caption "- Walls of Data $en Sensitising Participants to the World of Commercially-held Data and GDPR" "- Sentence Ranking $en Bringing Support Workers and Families to a Shared Problem Space"

# 3.3 (was 3.2's text, the quoted "Family Facts" caption)
Replace-Caption "- $($ldq)Family Facts$($rdq) $en What is Data?" "- Walls of Data $en Sensitising Participants to the World of Commercially-held Data and GDPR"

# 3.2 (was the dropped "Research Activities and Contexts" caption)
Replace-Caption "- Research Activities and Contexts" "- $($ldq)Family Facts$($rdq) $en What is Data?"
